# Apply updated price/volume figures to the symbol list (cryptos sheet).
# Source data refreshed on Fri Jan 13 20:49:58 UTC 2023.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ordered list of cell -> new text value (column D = Price, column E = Volume(1h)).
$updates = [ordered]@{
    "D2" = "288.26"
    "E2" = "0.76%"
    "D3" = "29.15"
    "E3" = "1.46%"
    "D4" = "5.285"
    "E4" = "4.78%"
    "D5" = "0.07061"
    "E5" = "5.48%"
    "D6" = "7.445"
    "E6" = "1.35%"
    "D7" = "3.558"
    "E7" = "5.08%"
    "D8" = "1.399"
    "E8" = "2.14%"
    "D9" = "0.9032"
    "E9" = "-3.96%"
    "D10" = "0.1602"
    "E10" = "2.17%"
    "D11" = "0.07598"
    "E11" = "12.95%"
    "D12" = "0.07670"
    "E12" = "1.04%"
    "D13" = "0.02911"
    "E13" = "-1.21%"
    "D14" = "0.09030"
    "E14" = "0.33%"
    "D15" = "0.001593"
    "E15" = "-0.37%"
    "D16" = "0.0006498"
    "E16" = "0.67%"
    "D17" = "0.006225"
    "E17" = "-6.52%"
    "E18" = "1.06%"
    "E19" = "-0.63%"
    "D20" = "0.3211"
    "E20" = "0.08%"
    "D21" = "0.1347"
    "E21" = "2.79%"
    "D22" = "4.003"
    "E22" = "-1.80%"
    "D23" = "0.1597"
    "E23" = "3.05%"
    "D24" = "0.04522"
    "E24" = "0.70%"
    "D25" = "0.001209"
    "E25" = "2.60%"
    "D26" = "0.004163"
    "E26" = "-7.31%"
    "D27" = "0.0001168"
    "E27" = "-6.34%"
    "D28" = "0.0001666"
    "E28" = "3.21%"
    "D40" = "0.04362"
    "E40" = "3.98%"
    "D41" = "0.006940"
    "E41" = "3.05%"
    "D42" = "0.1252"
    "E42" = "-0.02%"
    "D43" = "0.002066"
    "E43" = "2.55%"
    "D44" = "0.01169"
    "E44" = "-4.84%"
    "D45" = "0.00005866"
    "E45" = "3.69%"
    "D47" = "0.01297"
    "E47" = "-0.49%"
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    # Force text interpretation so values like "288.26" or "0.76%" are not
    # auto-converted to numbers/percentages by Excel's input parser, matching
    # the original inline-string cell type.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    # Drop the temporary text format again so the cell keeps its original
    # (default / General) style, exactly as before the edit.
    $cell.ClearFormats()
}
